$wb = $excel.ActiveWorkbook

# Rename "TFP" sheet to "FPROD"
$wsFPROD = $wb.Worksheets.Item("TFP")
$wsFPROD.Name = "FPROD"

# Add the new "factor specific productivity" column (GOS) next to the
# existing column, and rename the existing header from TFP to COE.
$wsFPROD.Range("B1").Value = "COE"
$wsFPROD.Range("C1").Value = "GOS"

# Fill C2:C36 with 1, mirroring the existing B2:B36 values.
for ($r = 2; $r -le 36; $r++) {
    $wsFPROD.Cells.Item($r, 3).Value = 1
}

# Update selection on the elasTRADE sheet (no longer the active tab).
$wsTRADE = $wb.Worksheets.Item("elasTRADE")
$wsTRADE.Select() | Out-Null
$wsTRADE.Range("E11").Select() | Out-Null

# Make FPROD the active sheet/tab with its own selection.
$wsFPROD.Select() | Out-Null
$wsFPROD.Range("B7").Select() | Out-Null
